$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Chinese/utterance text in column K (rows 2-13) ---
$ws.Range("K2").Value = "明白"
$ws.Range("K3").Value = "明白"
$ws.Range("K4").Value = "哦"
$ws.Range("K5").Value = "我諗下"
$ws.Range("K6").Value = "俾我諗下"
$ws.Range("K7").Value = "嗯，等我諗下"
$ws.Range("K8").Value = "嗯"
$ws.Range("K9").Value = "嗯"
$ws.Range("K10").Value = "好"
$ws.Range("K11").Value = "好啊"
$ws.Range("K12").Value = "好嘞"
$ws.Range("K13").Value = "我聽到"

# K3, K4, K8 and K10 need to pick up the "charset" CJK font variant
# (reuses the existing style rather than creating a new one).
$ws.Range("K3").Font.Name = "Noto Sans CJK SC"
$ws.Range("K3").Font.Size = 12
$ws.Range("K4").Font.Name = "Noto Sans CJK SC"
$ws.Range("K4").Font.Size = 12
$ws.Range("K8").Font.Name = "Noto Sans CJK SC"
$ws.Range("K8").Font.Size = 12
$ws.Range("K10").Font.Name = "Noto Sans CJK SC"
$ws.Range("K10").Font.Size = 12

# --- Row heights picked up new values after the content edit ---
$ws.Rows(2).RowHeight = 15
$ws.Rows(5).RowHeight = 15
$ws.Rows(6).RowHeight = 15
$ws.Rows(7).RowHeight = 15
$ws.Rows(8).RowHeight = 15
$ws.Rows(13).RowHeight = 15

# --- Move the active selection ---
$ws.Range("J22").Select()
